$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the existing header cell (E1) onto the new
# header cell F1, then set its text. xlPasteFormats = -4122
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

# Add time_taken values for each data row
$times = @(
    "2021-10-05 13:41:11.138411",
    "2021-10-05 13:41:11.138422",
    "2021-10-05 13:41:11.138426",
    "2021-10-05 13:41:11.138428",
    "2021-10-05 13:41:11.138432",
    "2021-10-05 13:41:11.138434",
    "2021-10-05 13:41:11.138437",
    "2021-10-05 13:41:11.138440",
    "2021-10-05 13:41:11.138443",
    "2021-10-05 13:41:11.138446",
    "2021-10-05 13:41:11.138449",
    "2021-10-05 13:41:11.138451"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
